$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.875.49'
$ws.Cells.Item(2, 5).Value = '  -0.99%  '

$ws.Cells.Item(3, 4).Value = '2.604.65'
$ws.Cells.Item(3, 5).Value = '  -1.36%  '

$ws.Cells.Item(4, 5).Value = '  +0.07%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '556.41'
$ws.Cells.Item(5, 5).Value = '  +3.61%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '143.72'
$ws.Cells.Item(6, 5).Value = '  -1.20%  '

$ws.Cells.Item(7, 5).Value = '  -0.07%  '

$ws.Cells.Item(8, 5).Value = '  +4.21%  '

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '6.81'
$ws.Cells.Item(9, 5).Value = '  -2.75%  '

$ws.Cells.Item(10, 5).Value = '  -0.97%  '

$ws.Cells.Item(11, 5).Value = '  +5.84%  '

$ws.Cells.Item(12, 5).Value = '  -1.10%  '

$ws.Cells.Item(13, 4).Value = '3.068.29'
$ws.Cells.Item(13, 5).Value = '  -1.35%  '

$ws.Cells.Item(14, 4).Value = '58.854.56'
$ws.Cells.Item(14, 5).Value = '  -0.91%  '

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '20.97'
$ws.Cells.Item(15, 5).Value = '  -1.56%  '

$ws.Cells.Item(16, 4).Value = '2.616.91'
$ws.Cells.Item(16, 5).Value = '  -2.01%  '

$ws.Cells.Item(17, 5).Value = '  -1.97%  '

$ws.Cells.Item(18, 5).Value = '  -0.85%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '337.61'
$ws.Cells.Item(19, 5).Value = '  -0.28%  '

$ws.Cells.Item(20, 5).Value = '  -2.03%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.18'
$ws.Cells.Item(21, 5).Value = '  -0.97%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.997'
$ws.Cells.Item(22, 5).Value = '  -0.34%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '66.29'
$ws.Cells.Item(23, 5).Value = '  -0.04%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '0.428'
$ws.Cells.Item(24, 5).Value = '  +2.57%  '

$ws.Cells.Item(25, 5).Value = '  -0.22%  '

$ws.Cells.Item(26, 5).Value = '  -2.20%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '7.15'
$ws.Cells.Item(27, 5).Value = '  -2.15%  '

$ws.Cells.Item(28, 4).Value = '0.0₃0761'
$ws.Cells.Item(28, 5).Value = '  +0.90%  '

$ws.Cells.Item(29, 5).Value = '  -0.05%  '

$ws.Cells.Item(30, 5).Value = '  +1.09%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '5.99'
$ws.Cells.Item(31, 5).Value = '  +1.37%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '154.02'
$ws.Cells.Item(32, 5).Value = '  +2.04%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '18.93'
$ws.Cells.Item(33, 5).Value = '  +0.47%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '3.96'
$ws.Cells.Item(34, 5).Value = '  -1.23%  '

$ws.Cells.Item(35, 2).Value = 'SuiNetwork'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.906'
$ws.Cells.Item(35, 5).Value = '  +6.26%  '

$ws.Cells.Item(36, 2).Value = 'Fetch.AI'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.907'
$ws.Cells.Item(36, 5).Value = '  +8.05%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.14'
$ws.Cells.Item(37, 5).Value = '  -0.57%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '37.12'
$ws.Cells.Item(38, 5).Value = '  -0.72%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.47'
$ws.Cells.Item(39, 5).Value = '  +0.98%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '3.59'
$ws.Cells.Item(40, 5).Value = '  -0.79%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '283.25'
$ws.Cells.Item(41, 5).Value = '  -0.95%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.998'
$ws.Cells.Item(42, 5).Value = '  -0.13%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.601'
$ws.Cells.Item(43, 5).Value = '  -0.06%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.0538'
$ws.Cells.Item(44, 5).Value = '  -0.14%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0954'
$ws.Cells.Item(45, 5).Value = '  +1.19%  '

$ws.Cells.Item(46, 5).Value = '  -1.30%  '

$ws.Cells.Item(47, 2).Value = 'RenderToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '4.67'
$ws.Cells.Item(47, 5).Value = '  +2.33%  '

$ws.Cells.Item(48, 2).Value = 'VeChain'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '0.0227'
$ws.Cells.Item(48, 5).Value = '  -0.14%  '

$ws.Cells.Item(49, 2).Value = 'Maker'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Cells.Item(49, 4).Value = '1.946.99'
$ws.Cells.Item(49, 5).Value = '  -0.84%  '

$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '118.16'
$ws.Cells.Item(50, 5).Value = '  +6.11%  '

$ws.Cells.Item(51, 5).Value = '  -2.71%  '

